$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tiny float re-write of the previous last row's timestamp (matches source diff).
$ws.Range("A88").Value = 44401.77015799537

# Append a new data row (row 89) with the latest retrieved data point.
$ws.Range("A89").Value = 44402.76951802956
$ws.Range("A89").NumberFormat = $ws.Range("A88").NumberFormat

$ws.Range("B89").Value = 79572
$ws.Range("C89").Value = 67135
$ws.Range("D89").Value = 3744
$ws.Range("E89").Value = 2186
$ws.Range("F89").Value = 1575
$ws.Range("G89").Value = 20839
$ws.Range("H89").Value = 1619
$ws.Range("I89").Value = 882
$ws.Range("J89").Value = 199
